$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 values
$ws.Range("B2").Value = 115
$ws.Range("C2").Value = "TVM"
$ws.Range("D2").Value = "BGLR"
$ws.Range("E2").Value = 1245

# Update row 3 values
$ws.Range("B3").Value = 116
$ws.Range("C3").Value = "Chennai"
$ws.Range("D3").Value = "Cochin"
$ws.Range("E3").Value = 6734

# Move the selection to E3
$ws.Range("E3").Select()
